$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header values for columns Q1:AB1 ---
# Q reuses the existing "Keyword" string (shared string index 2)
$newHeaders = @(
    "Keyword",
    "NewDescription",
    "Cat",
    "PartsClassID",
    "CurrentDescription",
    "Group",
    "U/I",
    "AW_StockStatus",
    "NewMfgNo",
    "Bin1",
    "Bin2",
    "Bin3"
)

$startCol = 17  # Q
for ($i = 0; $i -lt $newHeaders.Length; $i++) {
    $ws.Cells.Item(1, $startCol + $i).Value = $newHeaders[$i]
}

# Copy the existing header style (bold + centered) from A1 onto the new header
# cells so they match the rest of the row without creating extra style records.
$ws.Range("A1").Copy()
$ws.Range("Q1:AB1").PasteSpecial(-4122)

# --- Column width changes ---
$ws.Columns.Item(1).ColumnWidth = 10.877604166666666    # A
$ws.Columns.Item(7).ColumnWidth = 30.877604166666668    # G
$ws.Columns.Item(8).ColumnWidth = 10.877604166666666    # H
$ws.Columns.Item(18).ColumnWidth = 14.592447916666666   # R
$ws.Columns.Item(20).ColumnWidth = 10.877604166666666   # T
$ws.Columns.Item(21).ColumnWidth = 17.307291666666668   # U
$ws.Columns.Item(24).ColumnWidth = 14.736979166666666   # X
$ws.Columns.Item(25).ColumnWidth = 10.451822916666666   # Y
